# Auto-generated edit script: applies scheduled market-data refresh values
# to the Phantom_Profits leve-profit sheets (per-row currentAveragePrice /
# LevePrice / LeveProfit recalculation).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 9
$ws.Range("H9").Value = 1375
$ws.Range("I9").Value = 1375
$ws.Range("K9").Value = 1375
$ws.Range("M9").Value = -1206

# row 12
$ws.Range("H12").Value = 892.61536
$ws.Range("I12").Value = 892.61536
$ws.Range("K12").Value = 892.61536
$ws.Range("M12").Value = -722.61536

# row 74
$ws.Range("H74").Value = 4984.143
$ws.Range("I74").Value = 4984.143
$ws.Range("K74").Value = 4984.143
$ws.Range("M74").Value = -4048.143

# row 77
$ws.Range("H77").Value = 4984.143
$ws.Range("I77").Value = 4984.143
$ws.Range("K77").Value = 24920.715
$ws.Range("M77").Value = -20240.715

# row 112
$ws.Range("H112").Value = 3999.6365
$ws.Range("J112").Value = 3999.6365
$ws.Range("L112").Value = 11998.9095
$ws.Range("N112").Value = -14214.9095

# row 116
$ws.Range("H116").Value = 7000
$ws.Range("J116").Value = 7000
$ws.Range("L116").Value = 7000
$ws.Range("N116").Value = -13884

# row 137
$ws.Range("H137").Value = 3721.5789
$ws.Range("I137").Value = 3664.6428
$ws.Range("J137").Value = 3881
$ws.Range("K137").Value = 10993.9284
$ws.Range("L137").Value = 11643
$ws.Range("M137").Value = -8443.928400000001
$ws.Range("N137").Value = -16743

$ws = $wb.Worksheets.Item("ARM")
# row 5
$ws.Range("H5").Value = 277.2857
$ws.Range("I5").Value = 277.2857
$ws.Range("K5").Value = 277.2857
$ws.Range("M5").Value = -165.2857

# row 32
$ws.Range("H32").Value = 10893
$ws.Range("I32").Value = 10931.869
$ws.Range("K32").Value = 10931.869
$ws.Range("M32").Value = -10644.869

# row 97
$ws.Range("H97").Value = 1106.8462
$ws.Range("I97").Value = 308.0909
$ws.Range("J97").Value = 5500
$ws.Range("K97").Value = 308.0909
$ws.Range("L97").Value = 5500
$ws.Range("M97").Value = 187.9091
$ws.Range("N97").Value = -6492

# row 122
$ws.Range("H122").Value = 2065.889
$ws.Range("I122").Value = 2199.1875
$ws.Range("K122").Value = 6597.5625
$ws.Range("M122").Value = -4147.5625

# row 132
$ws.Range("H132").Value = 2744.3333
$ws.Range("I132").Value = 2744.3333
$ws.Range("K132").Value = 8232.999899999999
$ws.Range("M132").Value = -5702.999899999999

$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 277.2857
$ws.Range("I4").Value = 277.2857
$ws.Range("K4").Value = 277.2857
$ws.Range("M4").Value = -162.2857

# row 64
$ws.Range("H64").Value = 3226.8
$ws.Range("I64").Value = 1478.3334
$ws.Range("J64").Value = 3976.1428
$ws.Range("K64").Value = 1478.3334
$ws.Range("L64").Value = 3976.1428
$ws.Range("M64").Value = -1253.3334
$ws.Range("N64").Value = -4426.1428

# row 67
$ws.Range("H67").Value = 3226.8
$ws.Range("I67").Value = 1478.3334
$ws.Range("J67").Value = 3976.1428
$ws.Range("K67").Value = 1478.3334
$ws.Range("L67").Value = 3976.1428
$ws.Range("M67").Value = -698.3334
$ws.Range("N67").Value = -5536.1428

# row 107
$ws.Range("H107").Value = 1078.4166
$ws.Range("I107").Value = 1130.25
$ws.Range("J107").Value = 974.75
$ws.Range("K107").Value = 1130.25
$ws.Range("L107").Value = 974.75
$ws.Range("M107").Value = 789.75
$ws.Range("N107").Value = -4814.75

$ws = $wb.Worksheets.Item("CRP")
# row 3
$ws.Range("H3").Value = 1995
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1995
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1995
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -2221

# row 31
$ws.Range("H31").Value = 5937.1665
$ws.Range("I31").Value = 2124.6
$ws.Range("J31").Value = 25000
$ws.Range("K31").Value = 2124.6
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = -1829.6
$ws.Range("N31").Value = -25590

# row 34
$ws.Range("H34").Value = 5937.1665
$ws.Range("I34").Value = 2124.6
$ws.Range("J34").Value = 25000
$ws.Range("K34").Value = 2124.6
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = -1922.6
$ws.Range("N34").Value = -25404

# row 99
$ws.Range("H99").Value = 2979
$ws.Range("I99").Value = 2979
$ws.Range("K99").Value = 2979
$ws.Range("M99").Value = -1481

# row 105
$ws.Range("H105").Value = 3239.8
$ws.Range("I105").Value = 3047
$ws.Range("K105").Value = 3047
$ws.Range("M105").Value = -1300

# row 121
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# row 126
$ws.Range("H126").Value = 2979
$ws.Range("I126").Value = 2979
$ws.Range("K126").Value = 8937
$ws.Range("M126").Value = -6467

# row 132
$ws.Range("H132").Value = 3814.4285
$ws.Range("I132").Value = 3814.4285
$ws.Range("K132").Value = 11443.2855
$ws.Range("M132").Value = -8913.2855

$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 138.5
$ws.Range("I2").Value = 106.666664
$ws.Range("J2").Value = 170.33333
$ws.Range("K2").Value = 639.999984
$ws.Range("L2").Value = 1021.99998
$ws.Range("M2").Value = -526.999984
$ws.Range("N2").Value = -1247.99998

# row 4
$ws.Range("H4").Value = 3855569.5
$ws.Range("I4").Value = 120650.71
$ws.Range("K4").Value = 361952.13
$ws.Range("M4").Value = -361840.13

# row 15
$ws.Range("H15").Value = 192
$ws.Range("J15").Value = 98
$ws.Range("L15").Value = 294
$ws.Range("N15").Value = -574

# row 47
$ws.Range("H47").Value = 303.75
$ws.Range("I47").Value = 303.75
$ws.Range("K47").Value = 911.25
$ws.Range("M47").Value = -480.25

# row 51
$ws.Range("H51").Value = 764
$ws.Range("I51").Value = 764
$ws.Range("K51").Value = 2292
$ws.Range("M51").Value = -1832

# row 55
$ws.Range("H55").Value = 13600.333
$ws.Range("J55").Value = 19998.5
$ws.Range("L55").Value = 59995.5
$ws.Range("N55").Value = -60349.5

# row 92
$ws.Range("H92").Value = 246.33333
$ws.Range("I92").Value = 215.6
$ws.Range("J92").Value = 400
$ws.Range("K92").Value = 646.8
$ws.Range("L92").Value = 1200
$ws.Range("M92").Value = 601.2
$ws.Range("N92").Value = -3696

# row 97
$ws.Range("H97").Value = 1609.25
$ws.Range("I97").Value = 1479
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 4437
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -3941
$ws.Range("N97").Value = -6992

# row 122
$ws.Range("H122").Value = 535.9
$ws.Range("I122").Value = 494.83334
$ws.Range("K122").Value = 4453.50006
$ws.Range("M122").Value = -2003.50006

# row 131
$ws.Range("H131").Value = 1767.2858
$ws.Range("J131").Value = 2033
$ws.Range("L131").Value = 6099
$ws.Range("N131").Value = -16179

$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 133.28572
$ws.Range("I2").Value = 133.28572
$ws.Range("K2").Value = 133.28572
$ws.Range("M2").Value = -20.28572

# row 64
$ws.Range("H64").Value = 190271
$ws.Range("J64").Value = 190271
$ws.Range("L64").Value = 190271
$ws.Range("N64").Value = -190767

# row 67
$ws.Range("H67").Value = 190271
$ws.Range("J67").Value = 190271
$ws.Range("L67").Value = 190271
$ws.Range("N67").Value = -191987

# row 80
$ws.Range("H80").Value = 2391
$ws.Range("I80").Value = 2401.3333
$ws.Range("J80").Value = 2375.5
$ws.Range("K80").Value = 2401.3333
$ws.Range("L80").Value = 2375.5
$ws.Range("M80").Value = -1403.3333
$ws.Range("N80").Value = -4371.5

# row 83
$ws.Range("H83").Value = 2391
$ws.Range("I83").Value = 2401.3333
$ws.Range("J83").Value = 2375.5
$ws.Range("K83").Value = 12006.6665
$ws.Range("L83").Value = 11877.5
$ws.Range("M83").Value = -7014.666499999999
$ws.Range("N83").Value = -21861.5

# row 122
$ws.Range("H122").Value = 1713.2609
$ws.Range("I122").Value = 1455.0526
$ws.Range("K122").Value = 4365.1578
$ws.Range("M122").Value = -1915.1578

# row 132
$ws.Range("H132").Value = 4327.5835
$ws.Range("I132").Value = 3190.4
$ws.Range("K132").Value = 9571.200000000001
$ws.Range("M132").Value = -7041.200000000001

$ws = $wb.Worksheets.Item("LTW")
# row 100
$ws.Range("H100").Value = 1229
$ws.Range("I100").Value = 1229
$ws.Range("K100").Value = 1229
$ws.Range("M100").Value = -688

$ws = $wb.Worksheets.Item("WVR")
# row 9
$ws.Range("H9").Value = 224.66667
$ws.Range("I9").Value = 224.66667
$ws.Range("K9").Value = 224.66667
$ws.Range("M9").Value = -84.66667000000001

# row 122
$ws.Range("H122").Value = 3978.2
$ws.Range("I122").Value = 3978.2
$ws.Range("K122").Value = 11934.6
$ws.Range("M122").Value = -9484.599999999999

# row 136
$ws.Range("H136").Value = 5797.375
$ws.Range("I136").Value = 5345.625
$ws.Range("K136").Value = 16036.875
$ws.Range("M136").Value = -13486.875
